$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.005") are
# preserved verbatim instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.463.82"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "1.474.04"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("D5").Value = "0.9774"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").Value = "275.26"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "0.3650"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D8").Value = "0.3070"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "39.78"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "1.052"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "0.06642"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "0.9995"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "5.465"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "18.03"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "6.172"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "0.00001030"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "1.472.66"
$ws.Range("E17").Value = "  +3.52%  "
$ws.Range("D18").Value = "0.9845"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "0.05872"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "69.42"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "5.461"
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("D22").Value = "14.42"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "2.247"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "20.500.53"
$ws.Range("E25").Value = "  +2.41%  "
$ws.Range("D26").Value = "141.83"
$ws.Range("E26").Value = "  +6.38%  "
$ws.Range("D27").Value = "2.150"
$ws.Range("E27").Value = "  -6.06%  "
$ws.Range("D28").Value = "17.24"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "1.628.74"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("D30").Value = "113.84"
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("D31").Value = "3.852"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "4.982"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("D33").Value = "0.8017"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").Value = "0.07863"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "1.547"
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("D36").Value = "0.05762"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "4.749"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.157"
$ws.Range("E38").Value = "  +4.80%  "
$ws.Range("D39").Value = "7.782"
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").Value = "0.9760"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "0.1876"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "0.5295"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "3.492"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "12.01"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").Value = "117.71"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "0.5193"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "1.774"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "0.06450"
$ws.Range("E50").Value = "  +3.75%  "
$ws.Range("D51").Value = "0.9913"
$ws.Range("E51").Value = "  -0.93%  "
